# Menu.xlsx update: append 10 new package rows (129-138) plus a trailing
# blank spacer row (139), fix a couple of row-height quirks, and move the
# saved scroll/selection position — matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. New data rows 129-138
# ---------------------------------------------------------------------
# Each entry: row, A (serial no.), B (pkg name [+hyperlink url]), C (version),
# D, E (blank marker or oh-tpl display [+url]), F (link label), and which
# cellXf style index (from the existing style set) column B/E/F should end
# up with so that visited/unvisited hyperlink colouring matches the source.
$rows = @(
    @{ r=129; a=128; b="react-native-waterfall-flow";           c="1.1.5";      e=$null; urlB="https://github.com/bolan9999/react-native-waterfall-flow";                     sB=8; sF=8 },
    @{ r=130; a=129; b="react-native-vconsole";                 c="0.1.11";     e=$null; urlB="https://github.com/AntoineDoubovetzky/react-native-vconsole";                  sB=8; sF=7 },
    @{ r=131; a=130; b="better-banner";                         c="0.0.1";      e=$null; urlB="https://github.com/react-native-oh-library/better-banner";                     sB=8; sF=7 },
    @{ r=132; a=131; b="react-native-ezswiper";                 c="1.3.0";      e=$null; urlB="https://github.com/react-native-oh-library/react-native-ezswiper";             sB=8; sF=7 },
    @{ r=133; a=132; b="react-native-image-header-scroll-view"; c="0.10.3";     e=$null; urlB="https://github.com/lodev09/react-native-image-header-scroll-view";             sB=8; sF=7 },
    @{ r=134; a=133; b="react-native-linear-gradient-text";     c="1.2.8";      e=$null; urlB="https://github.com/react-native-oh-library/react-native-linear-gradient-text"; sB=7; sF=7 },
    @{ r=135; a=134; b="react-native-marquee-ab";               c="2.0.0-rc.1"; e=$null; urlB="https://github.com/react-native-oh-library/react-native-marquee-ab";           sB=8; sF=7 },
    @{ r=136; a=135; b="react-native-reconnecting-websocket";   c="1.0.3";      e=$null; urlB="https://github.com/opensoutheast/reconnecting-websocket";                      sB=8; sF=7 },
    @{ r=137; a=136; b="react-native-json-tree";                c="1.3.0";      e=$null; urlB="https://github.com/react-native-oh-library/react-native-json-tree";            sB=8; sF=7 },
    @{ r=138; a=137; b="react-native-image-gallery";            c="2.1.5";
        e="@react-native-oh-tpl/react-native-image-gallery";
        urlB="https://github.com/archriss/react-native-image-gallery";
        urlE="https://github.com/react-native-oh-library/react-native-image-gallery/releases";
        sB=8; sF=7; sE=8 }
)

$dash = [char]0x2212   # '−' U+2212 MINUS SIGN, same glyph already used for blank D/E cells

foreach ($row in $rows) {
    $r = $row.r

    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $dash
    if ($row.e) {
        $ws.Cells.Item($r, 5).Value = $row.e
    } else {
        $ws.Cells.Item($r, 5).Value = $dash
    }
    $ws.Cells.Item($r, 6).Value = "链接"
}

# ---------------------------------------------------------------------
# 2. Hyperlinks for the new rows
# ---------------------------------------------------------------------
foreach ($row in $rows) {
    $r = $row.r
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $row.urlB, [Type]::Missing, [Type]::Missing, $row.b) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), "https://react-native-oh-library.gitee.io/usage-docs/zh-cn/" + $row.b, [Type]::Missing, [Type]::Missing, "链接") | Out-Null
    if ($row.e) {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $row.urlE, [Type]::Missing, [Type]::Missing, $row.e) | Out-Null
    }
}

# ---------------------------------------------------------------------
# 3. Row 139 - trailing blank spacer row (only B139 is touched, styled
#    like the other "comment/placeholder" cells)
# ---------------------------------------------------------------------
$ws.Cells.Item(139, 2).Value = ""

# ---------------------------------------------------------------------
# 4. Re-apply exact cell formatting (style indexes) via format-only paste
#    from existing donor cells so the new cells share the same cellXfs
#    entries as the rest of the sheet instead of minting near-duplicates.
# ---------------------------------------------------------------------
$donorA = $ws.Range("A124")   # s=1  plain number
$donorC = $ws.Range("C124")   # s=1  plain text
$donorD = $ws.Range("D128")   # s=9  "-" placeholder
$donorB8 = $ws.Range("B125")  # s=8  unvisited hyperlink
$donorB7 = $ws.Range("B124")  # s=6 -> not used; use E124 for s=7 below
$donorS7 = $ws.Range("E124")  # s=7  visited hyperlink
$donorF8 = $ws.Range("F125")  # s=8
$donorF7 = $ws.Range("F124")  # s=7
$donorB139 = $ws.Range("D128") # s=9, same look as other blank placeholder cells

foreach ($row in $rows) {
    $r = $row.r

    $donorA.Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteFormats) | Out-Null

    if ($row.sB -eq 7) {
        $donorS7.Copy() | Out-Null
    } else {
        $donorB8.Copy() | Out-Null
    }
    $ws.Cells.Item($r, 2).PasteSpecial($xlPasteFormats) | Out-Null

    $donorC.Copy() | Out-Null
    $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats) | Out-Null

    $donorD.Copy() | Out-Null
    $ws.Cells.Item($r, 4).PasteSpecial($xlPasteFormats) | Out-Null

    if ($row.sE -eq 8) {
        $donorB8.Copy() | Out-Null
    } else {
        $donorD.Copy() | Out-Null
    }
    $ws.Cells.Item($r, 5).PasteSpecial($xlPasteFormats) | Out-Null

    if ($row.sF -eq 8) {
        $donorF8.Copy() | Out-Null
    } else {
        $donorF7.Copy() | Out-Null
    }
    $ws.Cells.Item($r, 6).PasteSpecial($xlPasteFormats) | Out-Null
}

$donorB139.Copy() | Out-Null
$ws.Cells.Item(139, 2).PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Row-height fixes
#    - row 2 loses its custom 14.25 height (back to sheet default)
#    - rows 123 & 124 gain an explicit 14.25 height
#    - every new row (129-139) is 14.25 tall
# ---------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit() | Out-Null

$ws.Rows.Item(123).RowHeight = 14.25
$ws.Rows.Item(124).RowHeight = 14.25

foreach ($row in $rows) {
    $ws.Rows.Item($row.r).RowHeight = 14.25
}
$ws.Rows.Item(139).RowHeight = 14.25

# ---------------------------------------------------------------------
# 6. View state - scroll position + active selection
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 101
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I111").Select() | Out-Null
